# Registration Page Test Data 3
# Adds a new row of registration test data (row 10) to the "reg" worksheet,
# including a mailto hyperlink on the new Email cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reg")

# New test data columns: A=Email, B=Name, C=Username, D=Password, E=URL
$email    = "ak@gmail.com"
$name     = "bingi18"
$username = "RoyalEnfield"
$password = "Matrix@321456"
$url      = "https://www.google.com/"

$ws.Range("A10").Value = $email
$ws.Range("B10").Value = $name
$ws.Range("C10").Value = $username
$ws.Range("D10").Value = $password
$ws.Range("E10").Value = $url

# Match the font already used by the rest of the data rows.
$ws.Range("B10:E10").Font.Name = "Helvetica"

# Hyperlink the e-mail address, mirroring the existing hyperlinked cells
# (applies the workbook's "Hyperlink" cell style automatically).
$ws.Hyperlinks.Add($ws.Range("A10"), "mailto:$email")

# Selection ends up past the new data, on E14.
$ws.Range("E14").Select()
